$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: writes $value into $cell as literal text (never auto-converted to a
# number/date by Excel's input parser), without leaving any stray style on
# the target cell. We stage the text in a scratch cell that is explicitly
# formatted as Text ("@"), then copy/paste-special "Values" only (so just the
# value + its text type moves over — not the scratch cell's number format),
# then clear the scratch cell.
function Set-TextValue {
    param($cell, $value)
    $scratch = $ws.Range("ZZ1")
    $scratch.NumberFormat = "@"
    $scratch.Value = $value
    $scratch.Copy()
    $ws.Range($cell).PasteSpecial(-4163)
    $scratch.Clear()
}

Set-TextValue "D2" "245.59"
Set-TextValue "D3" "24.18"
Set-TextValue "D4" "5.253"
Set-TextValue "D5" "0.05782"
Set-TextValue "D6" "6.497"
Set-TextValue "D7" "3.144"
Set-TextValue "D8" "0.8184"
Set-TextValue "D9" "0.8492"
Set-TextValue "B10" "WazirX"
Set-TextValue "C10" "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
Set-TextValue "D10" "0.1360"
Set-TextValue "E10" "9WazirXWRX"
Set-TextValue "B11" "MandalaExchangeToken"
Set-TextValue "C11" "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
Set-TextValue "D11" "0.06952"
Set-TextValue "E11" "10MandalaExchangeTokenMDX"
Set-TextValue "B12" "LiechtensteinCryptoassetsExchange"
Set-TextValue "C12" "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
Set-TextValue "D12" "0.03149"
Set-TextValue "E12" "11LiechtensteinCryptoassetsExchangeLCX"
Set-TextValue "B13" "BitrueCoin"
Set-TextValue "C13" "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
Set-TextValue "D13" "0.02882"
Set-TextValue "E13" "12BitrueCoinBTR"
Set-TextValue "B14" "BitMartToken"
Set-TextValue "C14" "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
Set-TextValue "D14" "0.09387"
Set-TextValue "E14" "13BitMartTokenBMX"
Set-TextValue "B15" "MCDex"
Set-TextValue "C15" "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
Set-TextValue "D15" "3.757"
Set-TextValue "E15" "14MCDexMCB"
Set-TextValue "B16" "BitForexToken"
Set-TextValue "C16" "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
Set-TextValue "D16" "0.001513"
Set-TextValue "E16" "15BitForexTokenBF"
Set-TextValue "B17" "CoinExToken"
Set-TextValue "C17" "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
Set-TextValue "D17" "0.04700"
Set-TextValue "E17" "16CoinExTokenCET"
Set-TextValue "B18" "One"
Set-TextValue "C18" "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
Set-TextValue "D18" "0.0005962"
Set-TextValue "E18" "17OneONE"
Set-TextValue "D19" "0.006276"
Set-TextValue "D20" "0.001238"
Set-TextValue "D21" "0.004610"
Set-TextValue "D22" "0.00006901"
Set-TextValue "E22" "21NitroExNTX"
Set-TextValue "D24" "2.145"
Set-TextValue "D26" "0.1347"
Set-TextValue "D40" "0.03652"
Set-TextValue "B41" "BKEXToken"
Set-TextValue "C41" "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue "D41" "0.1057"
Set-TextValue "E41" "40BKEXTokenBKK"
Set-TextValue "B42" "CEJI"
Set-TextValue "C42" "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextValue "D42" "0.002750"
Set-TextValue "E42" "41CEJICEJIBestin24h"
Set-TextValue "B43" "KickToken"
Set-TextValue "C43" "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
Set-TextValue "D43" "0.003019"
Set-TextValue "E43" "42KickTokenKICKWorstin24h"
Set-TextValue "D44" "0.007452"
Set-TextValue "D45" "0.00005266"
Set-TextValue "D47" "0.3613"
Set-TextValue "D48" "0.002329"

